# Apply the "cryptos list" price/volume update (GitHub Actions crypto scraper run).
# Cells D/E hold plain text (prices use dots as thousands separators, e.g. "42.769.86",
# and percentages are padded strings like "  -0.13%  "), so every write below targets
# Cells.Item(row, col).Value directly. Values that otherwise look like a plain decimal
# number (e.g. "301.49") are written with a leading apostrophe so Excel keeps them as
# text instead of silently converting them to a Number cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '42.769.86'
$ws.Cells.Item(2, 5).Value = '  -0.13%  '

$ws.Cells.Item(3, 4).Value = '2.313.43'
$ws.Cells.Item(3, 5).Value = '  +0.39%  '

$ws.Cells.Item(4, 5).Value = '  +0.00%  '

$ws.Cells.Item(5, 4).Value = '''301.49'
$ws.Cells.Item(5, 5).Value = '  -1.37%  '

$ws.Cells.Item(6, 4).Value = '''95.41'
$ws.Cells.Item(6, 5).Value = '  -1.50%  '

$ws.Cells.Item(7, 4).Value = '''0.503'
$ws.Cells.Item(7, 5).Value = '  -0.49%  '

$ws.Cells.Item(8, 5).Value = '  +0.00%  '

$ws.Cells.Item(9, 5).Value = '  -1.31%  '

$ws.Cells.Item(10, 5).Value = '  -3.51%  '

$ws.Cells.Item(11, 4).Value = '''18.91'
$ws.Cells.Item(11, 5).Value = '  +1.67%  '

$ws.Cells.Item(12, 5).Value = '  -0.24%  '

$ws.Cells.Item(13, 5).Value = '  +0.24%  '

$ws.Cells.Item(14, 5).Value = '  -1.92%  '

$ws.Cells.Item(15, 4).Value = '2.672.96'

$ws.Cells.Item(16, 4).Value = '2.275.92'
$ws.Cells.Item(16, 5).Value = '  -0.71%  '

$ws.Cells.Item(17, 4).Value = '''0.786'
$ws.Cells.Item(17, 5).Value = '  +0.81%  '

$ws.Cells.Item(18, 4).Value = '42.714.52'
$ws.Cells.Item(18, 5).Value = '  -0.06%  '

$ws.Cells.Item(19, 4).Value = '''12.10'
$ws.Cells.Item(19, 5).Value = '  -4.68%  '

$ws.Cells.Item(20, 5).Value = '  +1.69%  '

$ws.Cells.Item(21, 4).Value = '0.0₃0889'
$ws.Cells.Item(21, 5).Value = '  -0.49%  '

$ws.Cells.Item(22, 4).Value = '''67.67'
$ws.Cells.Item(22, 5).Value = '  +0.67%  '

$ws.Cells.Item(23, 4).Value = '''2.28'
$ws.Cells.Item(23, 5).Value = '  +5.47%  '

$ws.Cells.Item(24, 4).Value = '''234.99'
$ws.Cells.Item(24, 5).Value = '  -0.40%  '

$ws.Cells.Item(25, 5).Value = '  +0.14%  '

$ws.Cells.Item(26, 4).Value = '''2.41'
$ws.Cells.Item(26, 5).Value = '  -0.25%  '

$ws.Cells.Item(27, 4).Value = '''24.28'
$ws.Cells.Item(27, 5).Value = '  -1.80%  '

$ws.Cells.Item(28, 5).Value = '  +14.46%  '

$ws.Cells.Item(29, 2).Value = 'Cosmos'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(29, 4).Value = '''9.12'
$ws.Cells.Item(29, 5).Value = '  +0.84%  '

$ws.Cells.Item(30, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(30, 4).Value = '''32.15'
$ws.Cells.Item(30, 5).Value = '  -2.67%  '

$ws.Cells.Item(31, 2).Value = 'Monero'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(31, 4).Value = '''155.81'
$ws.Cells.Item(31, 5).Value = '  -6.14%  '

$ws.Cells.Item(33, 4).Value = '''4.99'
$ws.Cells.Item(33, 5).Value = '  +0.38%  '

$ws.Cells.Item(34, 4).Value = '''17.62'
$ws.Cells.Item(34, 5).Value = '  -2.21%  '

$ws.Cells.Item(35, 4).Value = '''4.44'
$ws.Cells.Item(35, 5).Value = '  +0.05%  '

$ws.Cells.Item(36, 5).Value = '  +1.69%  '

$ws.Cells.Item(37, 4).Value = '''2.32'
$ws.Cells.Item(37, 5).Value = '  -1.16%  '

$ws.Cells.Item(38, 5).Value = '  +2.34%  '

$ws.Cells.Item(39, 5).Value = '  -0.42%  '

$ws.Cells.Item(40, 5).Value = '  +0.38%  '

$ws.Cells.Item(41, 5).Value = '  -1.09%  '

$ws.Cells.Item(42, 4).Value = '''21.54'
$ws.Cells.Item(42, 5).Value = '  +19.37%  '

$ws.Cells.Item(43, 4).Value = '1.916.83'
$ws.Cells.Item(43, 5).Value = '  -4.15%  '

$ws.Cells.Item(44, 5).Value = '  -0.97%  '

$ws.Cells.Item(45, 4).Value = '''10.05'
$ws.Cells.Item(45, 5).Value = '  -1.78%  '

$ws.Cells.Item(46, 5).Value = '  -1.71%  '

$ws.Cells.Item(47, 4).Value = '''2.73'
$ws.Cells.Item(47, 5).Value = '  -1.15%  '

$ws.Cells.Item(48, 4).Value = '''2.87'
$ws.Cells.Item(48, 5).Value = '  +1.48%  '

$ws.Cells.Item(49, 4).Value = '2.545.27'
$ws.Cells.Item(49, 5).Value = '  +0.86%  '

$ws.Cells.Item(50, 4).Value = '''53.19'
$ws.Cells.Item(50, 5).Value = '  -0.71%  '

$ws.Cells.Item(51, 4).Value = '''72.17'
$ws.Cells.Item(51, 5).Value = '  +1.51%  '
